# "lowest point and reason analysis"
# The analyst re-opened the 可視化 (visualization) sheet, widened column A so the
# full survey-date labels (e.g. "1979年5月調査") are readable, and left the
# selection on A9 - the row flagged as the lowest point being investigated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("可視化")
$ws.Activate()

# Widen column A to fit the long Japanese date/survey labels.
$ws.Columns.Item(1).ColumnWidth = 21.86

# Move/leave the selection on A9 (the row under investigation).
$ws.Range("A9").Select()
